$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-04 Monday" "2025-08-05 Tuesday"

Replace-Text "882×2=1764" "775×2=1550"
Replace-Text "149×5=745" "957×5=4785"
Replace-Text "437×5=2185" "338×9=3042"
Replace-Text "641×4=2564" "180×3=540"
Replace-Text "919×4=3676" "276×6=1656"

Replace-Text "653×4=2612" "258×9=2322"
Replace-Text "916×2=1832" "832×2=1664"
Replace-Text "495×3=1485" "182×9=1638"
Replace-Text "321×3=963" "993×3=2979"
Replace-Text "213×7=1491" "102×6=612"

Replace-Text "666×3=1998" "494×3=1482"
Replace-Text "628×8=5024" "506×5=2530"
Replace-Text "652×7=4564" "828×9=7452"
Replace-Text "958×5=4790" "418×4=1672"
Replace-Text "442×8=3536" "504×9=4536"

Replace-Text "376×4=1504" "531×6=3186"
Replace-Text "973×9=8757" "910×8=7280"
Replace-Text "360×3=1080" "971×4=3884"
Replace-Text "490×5=2450" "353×2=706"
Replace-Text "233×5=1165" "403×5=2015"

Replace-Text "939×5=4695" "121×2=242"
Replace-Text "551×9=4959" "916×8=7328"
Replace-Text "932×2=1864" "814×7=5698"
Replace-Text "386×9=3474" "634×9=5706"
Replace-Text "111×3=333" "136×8=1088"
